# This script re-shuffles the data rows (2..102) of the sheet.
# Columns D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), R (Origen) and S (Precio $/Kg) are
# permuted across rows according to a fixed permutation derived from the
# target workbook, while all other columns remain untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 102

# Permutation: for target row i (array index 0 => row 2), the value is
# taken from this source row (1-based worksheet row number) in the
# ORIGINAL (pre-edit) data.
$sourceRows = @(13,14,28,17,30,99,4,7,73,40,21,20,41,88,2,101,89,84,36,85,82,8,83,77,18,96,65,26,27,68,90,97,63,95,16,79,23,60,56,69,19,92,58,39,64,47,11,10,100,44,45,54,57,32,50,9,74,51,15,53,34,37,22,46,78,49,24,71,33,86,87,61,102,94,42,43,6,29,31,72,38,93,76,81,66,70,75,48,62,55,25,80,59,12,67,98,52,5,35,3,91)

# Columns that move together with each row.
$cols = @(4, 13, 14, 15, 16, 18, 19)   # D, M, N, O, P, R, S

# 1) Snapshot all the original values for the columns that will change,
#    for every data row, BEFORE we overwrite anything.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Cells.Item($r, $c).Value2()
    }
    $snapshot[$r] = $rowVals
}

# 2) Write back the values according to the permutation.
for ($i = 0; $i -lt $sourceRows.Length; $i++) {
    $targetRow = $firstRow + $i
    $srcRow = $sourceRows[$i]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Cells.Item($targetRow, $c).Value = $srcVals[$c]
    }
}
